# Edit performed:
#  1. Slide 16's table switches to a different built-in table style
#     ({1ABF8A6F-88AD-4E92-83CC-07BBCBFD40CC} -> {B18ECC74-F4BE-4008-AAE7-33B8ADA17C40}).
#  2. The deck's colour theme (carried on the slide master / ppt/theme/theme1.xml)
#     is switched from the "Integral" palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{B18ECC74-F4BE-4008-AAE7-33B8ADA17C40}")

# --- 2. Swap the theme colour scheme (Integral -> Office) -----------------
# Slide.ThemeColorScheme edits the colour scheme carried by the presentation's
# slide master theme; changing it here applies to every slide in the deck.
$s1 = $p.Slides.Item(1)
$colors = $s1.ThemeColorScheme

$colors.Item(1).RGB  = 0x000000    # dk1      000000
$colors.Item(2).RGB  = 0xFFFFFF    # lt1      FFFFFF
$colors.Item(3).RGB  = 0x6A5444    # dk2      44546A (stored BGR)
$colors.Item(4).RGB  = 0xE6E6E7    # lt2      E7E6E6 (stored BGR)
$colors.Item(5).RGB  = 0xD59B5B    # accent1  5B9BD5 (stored BGR)
$colors.Item(6).RGB  = 0x317DED    # accent2  ED7D31 (stored BGR)
$colors.Item(7).RGB  = 0xA5A5A5    # accent3  A5A5A5
$colors.Item(8).RGB  = 0x00C0FF    # accent4  FFC000 (stored BGR)
$colors.Item(9).RGB  = 0xC47244    # accent5  4472C4 (stored BGR)
$colors.Item(10).RGB = 0x47AD70    # accent6  70AD47 (stored BGR)
$colors.Item(11).RGB = 0xC16305    # hlink    0563C1 (stored BGR)
$colors.Item(12).RGB = 0x724F95    # folHlink 954F72 (stored BGR)
